# Menus.xlsx edit: bump the resto_id values in column A (rows 2-100) from
# 1..99 up to 408..506 (offset +407), and (re)apply the alignment used on
# that column so its readingOrder is explicitly the context-dependent
# default (0), matching the style xf that Excel re-emits for the column
# after the values are rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$offset = 407
$firstRow = 2
$lastRow = 100

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $current = $cell.Value()
    $cell.Value = [double]($current + $offset)
}

# Re-assert the column's alignment (right/bottom, context reading order)
# so the style xf carries readingOrder="0" explicitly, same as the
# author's re-saved workbook.
$idRange = $ws.Range("A2:A100")
$idRange.HorizontalAlignment = -4152   # xlRight
$idRange.VerticalAlignment = -4107    # xlBottom
$idRange.ReadingOrder = 0              # xlContext (explicit, matches diff)
